$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle cells that lose the bold/red/bordered "peak" highlight but keep the gray data fill (style index 2) ---
$styleDonor2 = $ws.Range("B2")
foreach ($addr in @("D11","E15")) {
    $styleDonor2.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# --- Restyle cells that lose the bold/red/bordered "peak" highlight entirely, reverting to the plain default style (style index 0) ---
$styleDonor0 = $ws.Range("D16")
foreach ($addr in @("B16","C16","B17","C17","D17","E17","B18","C18","B19","C19","D19","E19","B20","C20","B21","C21","B22","C22","D22","E22","B23","C23","D23","B24","C24","D24","E24","B25","C25","E25","B26","C26","D26","E26","B27","C27","E27","B28","C28","D28","E28","B29","C29","D29")) {
    $styleDonor0.Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# --- Update numeric values to the recomputed results ---
$ws.Range("B2").Value = 0.000000000139724437396701196483074524341632359492049175742068
$ws.Range("C2").Value = 0.000000000139724003954060296244755318569822234026078788815539
$ws.Range("D2").Value = -0.000000004484665166410507400980220427122446835710434243083
$ws.Range("E2").Value = -0.00000000448466527228962580776776281066275942066567949950695
$ws.Range("B3").Value = 0.000000000051523545496083540458578014726564726737112476229186
$ws.Range("C3").Value = -0.000000000051522765870352403160582817071287990359773623083584
$ws.Range("D3").Value = -0.000000008061005447868465843622301026175080984614851331571117
$ws.Range("E3").Value = 0.000000008061006599490804792261557641143421282237113700830378
$ws.Range("B4").Value = -0.000000000034694469519536141888238489627838134765625
$ws.Range("C4").Value = -0.000000000034694469519536141888238489627838134765625
$ws.Range("D4").Value = 0.000000006834810495348619951982982456684112548828125
$ws.Range("E4").Value = 0.000000006834810495348619951982982456684112548828125
$ws.Range("B5").Value = -0.000000000869458648293357546000352492594027628558528419944196
$ws.Range("C5").Value = 0.000000000869427049682855805552284972776691182594355211676884
$ws.Range("D5").Value = -0.000000033137992893939560055155711193811129966491080267587677
$ws.Range("E5").Value = 0.000000033137958504206528087321820307334374611230032314779237
$ws.Range("B6").Value = -0.000000018767887200348951252790573470362323060811604591435753
$ws.Range("C6").Value = -0.000000018767903558672748410188321939451316922031764988787472
$ws.Range("D6").Value = -0.000003652095443572722841788120995420996450775419361889362335
$ws.Range("E6").Value = -0.000003652095439708135019940375698421775041424552910029888153
$ws.Range("B7").Value = -0.000000036535298066207322768679196557661814637185671017505229
$ws.Range("C7").Value = 0.000000036535265129923463191022981822256299189177752850810066
$ws.Range("D7").Value = -0.000006708911832358813302258377514464271484939672518521547318
$ws.Range("E7").Value = 0.000006708911841012671097548864668214108064603351522237062454
$ws.Range("B8").Value = 0.000000039996325129545183197338724520250585747760396770900115
$ws.Range("C8").Value = 0.000000000000897364047325970689836879704091160414591854532773
$ws.Range("D8").Value = -0.0000000002122941812300614958953440492979663871409012898539
$ws.Range("E8").Value = 0.000000000175141672310041495806137363951636692210866996788354
$ws.Range("B9").Value = -0.000000000001963188279622116083025166421073598639352592032381
$ws.Range("C9").Value = 0.000000000000475204508563751209853850965281196749856462702155
$ws.Range("D9").Value = -0.000000000154388885547582209998348362622877516514563467353582
$ws.Range("E9").Value = 0.000000000101925379540060695272220617643113994149617340667646
$ws.Range("B10").Value = -0.000000000000113121543496435299941028933846158879421598883463
$ws.Range("C10").Value = 0.00000000000000656590120850451086198084488243509326821367722
$ws.Range("D10").Value = -0.000000000002172551718702425047974506249220999671862425417501
$ws.Range("E10").Value = 0.000000000001319809822891875993299665763393960135054991389403
$ws.Range("B11").Value = -0.000000000212253071558845612303806544688452705854642132976551
$ws.Range("C11").Value = 0.000000000175107445440681604033048717875594834508801156403024
$ws.Range("D11").Value = 0.000039920168576844998974649525447233600061736069619655609131
$ws.Range("E11").Value = 0.000000034380794466569797589267586633565709774984497926197946
$ws.Range("B12").Value = 0.000000000000897364183035289936817982612681648957320135195914
$ws.Range("C12").Value = 0.000000039996325129784278099035952063616816687385835393797606
$ws.Range("D12").Value = 0.000000000175141672106477491086271567306636119665252948607304
$ws.Range("E12").Value = -0.000000000212294181530560702799373634309845326056942482750856
$ws.Range("B13").Value = 0.000000000000475206990105588868936875579507272959745023399591
$ws.Range("C13").Value = -0.00000000000196318910680272782825927434616620196770539474862
$ws.Range("D13").Value = 0.000000000101925377885699405542679912195009510680887565570174
$ws.Range("E13").Value = -0.000000000154388887201943409255009570083580060922967547298867
$ws.Range("B14").Value = 0.000000000000006565914133201582003038264748195906861383264186
$ws.Range("C14").Value = -0.000000000000113121524109389705851217287583706545619501547539
$ws.Range("D14").Value = 0.000000000001319809829354224932766977946484884912656601607672
$ws.Range("E14").Value = -0.000000000002172551712240076108507194066130074894260815199232
$ws.Range("B15").Value = 0.000000000175107570758544405816741460894011683108173471623559
$ws.Range("C15").Value = -0.00000000021225305997831708826020738105296069633820543742786
$ws.Range("D15").Value = 0.000000034380794352005282750673253663875605923294642707332969
$ws.Range("E15").Value = 0.000039920168576821417577397965725793937963317148387432098389
$ws.Range("B16").Value = 0.000000069862218698350624090931404452931019477546215057373047
$ws.Range("C16").Value = 0.000000069862001857916126990044602962370845489203929901123047
$ws.Range("D16").Value = -0.000002242332582702327888057869387239406933076679706573486328
$ws.Range("E16").Value = -0.000002242332636912436935849564714762216510735015617683529854
$ws.Range("B17").Value = 0.000000025761772748041771185716590627720612438622538320487365
$ws.Range("C17").Value = -0.000000025761382935176199150448359161125200245123778586275876
$ws.Range("D17").Value = -0.000004030502723934232895341370911390654896422347519546747208
$ws.Range("E17").Value = 0.000004030503299745403388747513884204920486808987334370613098
$ws.Range("B18").Value = -0.000000016812751781973808806300006734625596127585822614491917
$ws.Range("C18").Value = -0.000000016812735158952218940655852518796520289612317355931737
$ws.Range("D18").Value = 0.000003417644708441300858189210165605231850349809974431991577
$ws.Range("E18").Value = 0.000003417644704629653019061332266947950131452671485021710396
$ws.Range("B19").Value = -0.00000043472572156158122100050745367783022743424226064234972
$ws.Range("C19").Value = 0.000000434709922259577999173338564223301183631065214285627007
$ws.Range("D19").Value = -0.00001656833491067372920719884266116395110657322220504283905
$ws.Range("E19").Value = 0.000016568317715806359519950183289971334943402325734496116638
$ws.Range("B20").Value = -0.000009383943600470290302695275386213324964046478271484375
$ws.Range("C20").Value = -0.00000938395175909163825611614129229565151035785675048828125
$ws.Range("D20").Value = -0.001826047721776961896278379526847857050597667694091796875
$ws.Range("E20").Value = -0.00182604771985250304011660915648462832905352115631103515625
$ws.Range("B21").Value = -0.000018267649118836918421598589090670827772555639967322349548
$ws.Range("C21").Value = 0.000018267632650691870783900977448510616341081913560628890991
$ws.Range("D21").Value = -0.003354455919446954013235506764090132492128759622573852539062
$ws.Range("E21").Value = 0.003354455923773879983534884630103078961838036775588989257812
$ws.Range("B22").Value = 0.000019998162564772340268112044148196559945063199847936630249
$ws.Range("C22").Value = 0.000000000448682024322144895546633780458734008789178915321827
$ws.Range("D22").Value = -0.000000106147090615457302345496352381781779428138179355300963
$ws.Range("E22").Value = 0.000000087570836154154806123999302269744759996683569625020027
$ws.Range("B23").Value = -0.000000000981594295321012999533394349782611243426799774169922
$ws.Range("C23").Value = 0.00000000023760205906725099071360592842391736695395465517322
$ws.Range("D23").Value = -0.000000077194442939227227509779716285720496671274304389953613
$ws.Range("E23").Value = 0.00000005096268976010416736077690558204267290420830249786377
$ws.Range("B24").Value = -0.000000000056560770636693708011809476780284186396563494980683
$ws.Range("C24").Value = 0.000000000003282950811047408973144124632227764022402577426263
$ws.Range("D24").Value = -0.00000000108627586134161608458489050720464225152639414773148
$ws.Range("E24").Value = 0.000000000659904909843275660802537629923136819343199022114277
$ws.Range("B25").Value = -0.000000106126535671889299740037343206672248996369489759672433
$ws.Range("C25").Value = 0.000000087553722799750153746312086866510071558877825736999512
$ws.Range("D25").Value = 0.0199600842884224104201162930394275463186204433441162109375
$ws.Range("E25").Value = 0.000017190397233352730992994930603323666673531988635659217834
$ws.Range("B26").Value = 0.000000000448682090496593899788847770171429374386207200586796
$ws.Range("C26").Value = 0.000019998162564893110225731562290540921367210103198885917664
$ws.Range("D26").Value = 0.00000008757083605323877139252992593543289956414810148999095
$ws.Range("E26").Value = -0.000000106147090764349799370151533010192546058192419877741486
$ws.Range("B27").Value = 0.000000000237603329616671923862902815471898065879940986633301
$ws.Range("C27").Value = -0.000000000981594718837486626683563883943861583247780799865723
$ws.Range("D27").Value = 0.000000050962688913071220106476566513720172224566340446472168
$ws.Range("E27").Value = -0.000000077194443786260174764080055354042997350916266441345215
$ws.Range("B28").Value = 0.000000000003282957428492308993468740130339256211300380527973
$ws.Range("C28").Value = -0.000000000056560760710526357375477378323379881557009252190937
$ws.Range("D28").Value = 0.000000000659904913151998111014648329408771587623050436377525
$ws.Range("E28").Value = -0.000000001086275858032893013987320392948251329690023680996092
$ws.Range("B29").Value = 0.000000087553785268430013750962093155294496682472527027130127
$ws.Range("C29").Value = -0.000000106126530166175195526644342656347186348853028903249651
$ws.Range("D29").Value = 0.000017190397175966248816516082631444817252486245706677436829
$ws.Range("E29").Value = 0.0199600842884107600172516328029814758338034152984619140625
